$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the third (duplicate) data row entirely - the sheet shrinks from
# 3 data rows back to 2.
$ws.Rows("4:4").Delete()

# Replace the two remaining data rows with the newly imported records.
$ws.Range("A2").Value = "38895366700"
$ws.Range("B2").Value = "10061953"
$ws.Range("C2").Value = "SEI-490002/002287/2025"
$ws.Range("D2").Value = "ANTÔNIO"
$ws.Range("E2").Value = "D:\PROJETOS EM ANDAMENTO\consulta_receita\PDFs\RECEITA_38895366700_ANTÔNIO.pdf"

$ws.Range("A3").Value = "79348289772"
$ws.Range("B3").Value = "03041963"
$ws.Range("C3").Value = "SEI-490002/001100/2025"
$ws.Range("D3").Value = "ROSÂNGELA"
$ws.Range("E3").Value = "D:\PROJETOS EM ANDAMENTO\consulta_receita\PDFs\RECEITA_79348289772_ROSÂNGELA.pdf"

# Add the new UNIDADE column, matching the bold/centered header look already
# used for A1:E1, plus left/right borders around the header cell.
$ws.Range("F1").Value = "UNIDADE"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.Item(7).LineStyle = 1
$ws.Range("F1").Borders.Item(10).LineStyle = 1

$ws.Range("F2").Value = "SHF"
$ws.Range("F3").Value = "SHF"

# Columns go back to their default (unsized) widths.
$ws.Columns("A:E").ColumnWidth = 8.43

$ws.Range("F4").Select()
